$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "VIX Data"

# Scroll the active sheet view so row 353 is at the top and F371 is selected
$ws.Range("F371").Select()
$excel.ActiveWindow.ScrollRow = 353

# Set the window position (xWindow) for the workbook view
$excel.Left = 4240
